# Apply latest scraped crypto price / 1h volume-change values (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    # Writes $Text into $Address keeping it as plain text (matches the sheets
    # existing inline/shared string cells for Price and Volume(1h) columns).
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.Value = $Text
}

function Set-NumericLookingTextCellValue {
    # Some Price values (e.g. "48.51") would otherwise be auto-converted to a
    # number by Excel. Temporarily force text format, assign, then drop back to
    # the default (unstyled) cell style so no visible formatting changes stick.
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextCellValue "D2" "69.455.47"
Set-TextCellValue "E2" "  +0.31%  "
Set-TextCellValue "D3" "3.421.10"
Set-TextCellValue "E4" "  +0.09%  "
Set-NumericLookingTextCellValue "D5" "581.95"
Set-TextCellValue "E5" "  -0.80%  "
Set-NumericLookingTextCellValue "D6" "176.25"
Set-TextCellValue "E6" "  -2.31%  "
Set-TextCellValue "E7" "  +0.12%  "
Set-TextCellValue "D8" "3.412.99"
Set-TextCellValue "E8" "  +1.00%  "
Set-NumericLookingTextCellValue "D9" "0.592"
Set-TextCellValue "E9" "  -0.59%  "
Set-NumericLookingTextCellValue "D10" "0.199"
Set-TextCellValue "E10" "  +0.55%  "
Set-NumericLookingTextCellValue "D11" "0.581"
Set-TextCellValue "E11" "  -1.27%  "
Set-NumericLookingTextCellValue "D12" "48.51"
Set-TextCellValue "E12" "  -0.51%  "
Set-NumericLookingTextCellValue "D13" "0.0000280"
Set-TextCellValue "E13" "  -2.01%  "
Set-NumericLookingTextCellValue "D14" "695.92"
Set-TextCellValue "E14" "  +1.52%  "
Set-TextCellValue "D15" "3.966.96"
Set-TextCellValue "E15" "  +0.94%  "
Set-NumericLookingTextCellValue "D16" "8.60"
Set-TextCellValue "E16" "  -0.25%  "
Set-TextCellValue "D17" "69.471.44"
Set-TextCellValue "E17" "  +0.36%  "
Set-TextCellValue "D18" "3.417.87"
Set-TextCellValue "E18" "  +1.07%  "
Set-TextCellValue "E19" "  +0.85%  "
Set-NumericLookingTextCellValue "D20" "17.61"
Set-TextCellValue "E20" "  -0.57%  "
Set-NumericLookingTextCellValue "D21" "11.32"
Set-TextCellValue "E21" "  -0.42%  "
Set-NumericLookingTextCellValue "D22" "0.893"
Set-TextCellValue "E22" "  -0.91%  "
Set-NumericLookingTextCellValue "D23" "5.39"
Set-TextCellValue "E23" "  -0.76%  "
Set-NumericLookingTextCellValue "D24" "16.88"
Set-TextCellValue "E24" "  -1.21%  "
Set-NumericLookingTextCellValue "D25" "101.49"
Set-TextCellValue "E25" "  -2.48%  "
Set-TextCellValue "E26" "  -1.03%  "
Set-TextCellValue "E27" "  -2.67%  "
Set-NumericLookingTextCellValue "D28" "9.52"
Set-TextCellValue "E28" "  -0.91%  "
Set-NumericLookingTextCellValue "D29" "33.44"
Set-TextCellValue "E29" "  -2.64%  "
Set-NumericLookingTextCellValue "D30" "8.73"
Set-TextCellValue "E30" "  +0.40%  "
Set-NumericLookingTextCellValue "D31" "7.03"
Set-TextCellValue "E31" "  +0.74%  "
Set-NumericLookingTextCellValue "D32" "576.72"
Set-TextCellValue "E32" "  +3.64%  "
Set-NumericLookingTextCellValue "D33" "3.69"
Set-TextCellValue "E33" "  +0.33%  "
Set-NumericLookingTextCellValue "D34" "10.99"
Set-TextCellValue "E34" "  -1.87%  "
Set-NumericLookingTextCellValue "D35" "58.33"
Set-TextCellValue "E35" "  +0.50%  "
Set-TextCellValue "E36" "  -3.22%  "
Set-NumericLookingTextCellValue "D37" "1.00"
Set-TextCellValue "E37" "  +0.03%  "
Set-TextCellValue "D38" "3.534.66"
Set-TextCellValue "E38" "  -4.48%  "
Set-TextCellValue "E39" "  -1.59%  "
Set-NumericLookingTextCellValue "D40" "34.77"
Set-TextCellValue "E40" "  -0.19%  "
Set-TextCellValue "D41" "0.0₃0729"
Set-TextCellValue "E41" "  +3.60%  "
Set-NumericLookingTextCellValue "D42" "3.25"
Set-TextCellValue "E42" "  +0.51%  "
Set-TextCellValue "E43" "  -0.89%  "
Set-NumericLookingTextCellValue "D44" "0.331"
Set-TextCellValue "E44" "  -2.37%  "
Set-NumericLookingTextCellValue "D45" "0.0416"
Set-TextCellValue "E45" "  -0.93%  "
Set-NumericLookingTextCellValue "D46" "1.45"
Set-TextCellValue "E46" "  +4.36%  "
Set-NumericLookingTextCellValue "D47" "2.64"
Set-TextCellValue "E47" "  -0.37%  "
Set-TextCellValue "E48" "  -1.17%  "
Set-TextCellValue "E49" "  -0.30%  "
Set-NumericLookingTextCellValue "D50" "132.62"
Set-TextCellValue "E50" "  +0.25%  "
Set-TextCellValue "E51" "  +1.83%  "
